# Auto-generated edit script: refresh market-data derived values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 343.9
$ws.Range("I6").Value = 343.9
$ws.Range("K6").Value = 1031.7
$ws.Range("M6").Value = -919.6999999999998
$ws.Range("H29").Value = 2392.054
$ws.Range("J29").Value = 2841.9355
$ws.Range("L29").Value = 8525.806500000001
$ws.Range("N29").Value = -9087.806500000001
$ws.Range("H38").Value = 1215.2
$ws.Range("I38").Value = 239.11111
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 717.3333299999999
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = -345.3333299999999
$ws.Range("N38").Value = -30744
$ws.Range("H43").Value = 871.4286
$ws.Range("J43").Value = 966.6667
$ws.Range("L43").Value = 966.6667
$ws.Range("N43").Value = -1104.6667
$ws.Range("H58").Value = 3980.2
$ws.Range("I58").Value = 2300.3333
$ws.Range("J58").Value = 6500
$ws.Range("K58").Value = 6900.999899999999
$ws.Range("L58").Value = 19500
$ws.Range("M58").Value = -6750.999899999999
$ws.Range("N58").Value = -19800
$ws.Range("H87").Value = 32494.5
$ws.Range("J87").Value = 32494.5
$ws.Range("L87").Value = 32494.5
$ws.Range("N87").Value = -34990.5
$ws.Range("H90").Value = 32494.5
$ws.Range("J90").Value = 32494.5
$ws.Range("L90").Value = 97483.5
$ws.Range("N90").Value = -109963.5
$ws.Range("H129").Value = 929.24
$ws.Range("I129").Value = 334.66666
$ws.Range("J129").Value = 1117
$ws.Range("K129").Value = 1003.99998
$ws.Range("L129").Value = 3351
$ws.Range("M129").Value = 3996.00002
$ws.Range("N129").Value = -13351
$ws.Range("H132").Value = 2305.5334
$ws.Range("I132").Value = 1205.238
$ws.Range("J132").Value = 4872.8887
$ws.Range("K132").Value = 3615.714
$ws.Range("L132").Value = 14618.6661
$ws.Range("M132").Value = -1085.714
$ws.Range("N132").Value = -19678.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1317
$ws.Range("I2").Value = 978.2222
$ws.Range("K2").Value = 978.2222
$ws.Range("M2").Value = -865.2222
$ws.Range("H6").Value = 131666.89
$ws.Range("I6").Value = 165857.72
$ws.Range("J6").Value = 11999
$ws.Range("K6").Value = 165857.72
$ws.Range("L6").Value = 11999
$ws.Range("M6").Value = -165684.72
$ws.Range("N6").Value = -12345
$ws.Range("H32").Value = 4226.106
$ws.Range("I32").Value = 3032.5933
$ws.Range("J32").Value = 14285.714
$ws.Range("K32").Value = 3032.5933
$ws.Range("L32").Value = 14285.714
$ws.Range("M32").Value = -2745.5933
$ws.Range("N32").Value = -14859.714
$ws.Range("H50").Value = 1451.6666
$ws.Range("I50").Value = 1052.5
$ws.Range("J50").Value = 2250
$ws.Range("K50").Value = 1052.5
$ws.Range("L50").Value = 2250
$ws.Range("M50").Value = -338.5
$ws.Range("N50").Value = -3678
$ws.Range("H116").Value = 1317
$ws.Range("I116").Value = 978.2222
$ws.Range("K116").Value = 978.2222
$ws.Range("M116").Value = 1315.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1317
$ws.Range("I3").Value = 978.2222
$ws.Range("K3").Value = 978.2222
$ws.Range("M3").Value = -864.2222
$ws.Range("H14").Value = 1800
$ws.Range("I14").Value = 1800
$ws.Range("K14").Value = 1800
$ws.Range("M14").Value = -1628
$ws.Range("H134").Value = 4602.6743
$ws.Range("I134").Value = 4745.645
$ws.Range("J134").Value = 4233.3335
$ws.Range("K134").Value = 14236.935
$ws.Range("L134").Value = 12700.0005
$ws.Range("M134").Value = -11701.935
$ws.Range("N134").Value = -17770.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 655.4286
$ws.Range("I12").Value = 655.4286
$ws.Range("K12").Value = 655.4286
$ws.Range("M12").Value = -485.4286
$ws.Range("H99").Value = 1877.9131
$ws.Range("I99").Value = 1929.5294
$ws.Range("J99").Value = 1731.6666
$ws.Range("K99").Value = 1929.5294
$ws.Range("L99").Value = 1731.6666
$ws.Range("M99").Value = -431.5293999999999
$ws.Range("N99").Value = -4727.6666
$ws.Range("H126").Value = 1877.9131
$ws.Range("I126").Value = 1929.5294
$ws.Range("J126").Value = 1731.6666
$ws.Range("K126").Value = 5788.5882
$ws.Range("L126").Value = 5194.9998
$ws.Range("M126").Value = -3318.5882
$ws.Range("N126").Value = -10134.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2174.0833
$ws.Range("I17").Value = 158.14285
$ws.Range("J17").Value = 4996.4
$ws.Range("K17").Value = 474.42855
$ws.Range("L17").Value = 14989.2
$ws.Range("M17").Value = -305.42855
$ws.Range("N17").Value = -15327.2
$ws.Range("H68").Value = 1928.7407
$ws.Range("I68").Value = 1099.875
$ws.Range("J68").Value = 2277.7368
$ws.Range("K68").Value = 3299.625
$ws.Range("L68").Value = 6833.2104
$ws.Range("M68").Value = -2488.625
$ws.Range("N68").Value = -8455.2104
$ws.Range("H71").Value = 1928.7407
$ws.Range("I71").Value = 1099.875
$ws.Range("J71").Value = 2277.7368
$ws.Range("K71").Value = 9898.875
$ws.Range("L71").Value = 20499.6312
$ws.Range("M71").Value = -5842.875
$ws.Range("N71").Value = -28611.6312
$ws.Range("H106").Value = 4126.087
$ws.Range("J106").Value = 4126.087
$ws.Range("L106").Value = 12378.261
$ws.Range("N106").Value = -14270.261
$ws.Range("H112").Value = 111115500
$ws.Range("I112").Value = 1700
$ws.Range("J112").Value = 125004730
$ws.Range("K112").Value = 5100
$ws.Range("L112").Value = 375014190
$ws.Range("M112").Value = -3992
$ws.Range("N112").Value = -375016406
$ws.Range("H131").Value = 881.79
$ws.Range("I131").Value = 600.2
$ws.Range("J131").Value = 896.61053
$ws.Range("K131").Value = 1800.6
$ws.Range("L131").Value = 2689.83159
$ws.Range("M131").Value = 3239.4
$ws.Range("N131").Value = -12769.83159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4192.0933
$ws.Range("I132").Value = 2416.4546
$ws.Range("J132").Value = 6052.2856
$ws.Range("K132").Value = 7249.3638
$ws.Range("L132").Value = 18156.8568
$ws.Range("M132").Value = -4719.3638
$ws.Range("N132").Value = -23216.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 4999.5
$ws.Range("J12").Value = 4999.5
$ws.Range("L12").Value = 4999.5
$ws.Range("N12").Value = -5339.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 10005
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H41").Value = 24840.125
$ws.Range("J41").Value = 6820.857
$ws.Range("L41").Value = 6820.857
$ws.Range("N41").Value = -7600.857
